$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the remaining task row (row 3) ---
$ws.Range("C3").Value = "asd"
$ws.Range("D3").Value = 3

# F3: "17-11-2023" is not a valid date (day 17 can't be a month) so Excel
# keeps it as plain text automatically.
$ws.Range("F3").Value = "17-11-2023"

# G3: "12-12-2023" *is* a valid date, so a naive .Value assignment would be
# auto-converted into a date serial number. Build it as a text formula in a
# scratch cell, then copy/paste-special (values only) into G3 so it lands as
# a plain shared-string text value, matching the source workbook layout.
$ws.Range("ZZ1").Formula = "=""12-12-2023"""
$ws.Range("ZZ1").Copy()
$ws.Range("G3").PasteSpecial(-4163)
$ws.Range("ZZ1").ClearContents()

# --- Remove the other task rows (2, 2.1, 2.1.1) that were dropped ---
$ws.Range("A4:I6").Delete()
